# Fill in the 2x3 block of values that was added to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1
$ws.Range("B1").Value = 2
$ws.Range("C1").Value = 3
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 1

# Match the saved selection/active cell (G13) recorded in the sheet view.
$null = $ws.Range("G13").Select()
